# "added 4wk low sales check"
# Updates the forecast figures on "Forecast Comparison" (inventory coverage,
# seasonality index, a few MyForecast values, and the stockout
# risk/reorder urgency flags once coverage drops to 0) and refreshes the
# dependent roll-up totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison -----------------------------------------------

# Row 2 (W10)
$ws1.Range("D2").Value = 7
$ws1.Range("H2").Value = 10.52
$ws1.Range("L2").Value = 1.15

# Row 3 (W11)
$ws1.Range("H3").Value = 9.52
$ws1.Range("L3").Value = 1.05

# Row 4 (W12)
$ws1.Range("H4").Value = 8.52
$ws1.Range("L4").Value = 1.12

# Row 5 (W13)
$ws1.Range("H5").Value = 7.52
$ws1.Range("L5").Value = 1.1

# Row 6 (W14)
$ws1.Range("H6").Value = 6.52
$ws1.Range("L6").Value = 1.12

# Row 7 (W15)
$ws1.Range("H7").Value = 5.38
$ws1.Range("L7").Value = 0.87

# Row 8 (W16)
$ws1.Range("D8").Value = 7
$ws1.Range("H8").Value = 4.38
$ws1.Range("L8").Value = 1.13

# Row 9 (W17)
$ws1.Range("H9").Value = 3.38
$ws1.Range("L9").Value = 0.9

# Row 10 (W18)
$ws1.Range("H10").Value = 2.38
$ws1.Range("L10").Value = 1.09

# Row 11 (W19)
$ws1.Range("H11").Value = 1.38
$ws1.Range("L11").Value = 1.17

# Row 12 (W20) - coverage falls enough to flip risk/urgency flags
$ws1.Range("H12").Value = 0.38
$ws1.Range("I12").Value = "High"
$ws1.Range("J12").Value = "Urgent"
$ws1.Range("L12").Value = 1.2

# Row 13 (W21) - coverage hits 0, stockout risk flips to High
$ws1.Range("H13").Value = 0
$ws1.Range("I13").Value = "High"
$ws1.Range("L13").Value = 0.87

# Row 14 (W22)
$ws1.Range("D14").Value = 7
$ws1.Range("L14").Value = 0.97

# Row 15 (W23)
$ws1.Range("D15").Value = 7
$ws1.Range("L15").Value = 0.9

# Row 16 (W24)
$ws1.Range("D16").Value = 8
$ws1.Range("L16").Value = 0.82

# Row 17 (W25)
$ws1.Range("D17").Value = 8

# --- Summary -------------------------------------------------------------
# Roll-up totals recomputed after the MyForecast adjustments above.

$ws2.Range("B9").Value = "126"
$ws2.Range("B10").Value = "62"
$ws2.Range("B11").Value = "31"
$ws2.Range("B12").Value = "8"
$ws2.Range("B14").Value = "8"
